# Trade #76 closed at 2026-02-17 15:49:11 - unknown UNKNOWN +0.000%
#
# Applies:
#  - Summary sheet roll-up updates (Current Capital, Total P&L $/%% , Total
#    Trades, Losing Trades, Win Rate %%)
#  - Strategy Status roll-up updates for the MarketMaking row
#  - A new trade row (row 77 / trade #76) appended to both the "All Trades"
#    and "MarketMaking" logs, with the sheet dimension growing accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.82    # Current Capital
$summary.Range("B4").Value = -0.19      # Total P&L $
$summary.Range("B5").Value = -0.05      # Total P&L %
$summary.Range("B6").Value = 76         # Total Trades
$summary.Range("B8").Value = 41         # Losing Trades
$summary.Range("B9").Value = 31.58      # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.81999999999999   # Capital
$status.Range("D4").Value = 76                  # Trades
$status.Range("E4").Value = -0.19               # P&L $
$status.Range("F4").Value = -0.18               # P&L %
$status.Range("G4").Value = 31.58               # Win Rate %

# ---------------------------------------------------------------------
# Append new trade row (#76) to a trade-log sheet. Used for both
# "All Trades" and "MarketMaking", which carry identical logs.
# ---------------------------------------------------------------------
function Add-TradeRow {
    param($ws)

    $row = 77

    $ws.Cells.Item($row, 1).Value = 76              # Trade #
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"   # Date (force text, not a date serial)
    $ws.Cells.Item($row, 2).ClearFormats()
    $ws.Cells.Item($row, 3).Value = "15:49:05"      # Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"  # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"          # Side
    $ws.Cells.Item($row, 6).Value = 0.34            # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.05            # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"        # Status
    $ws.Cells.Item($row, 9).Value = -85.2941        # P&L %
    $ws.Cells.Item($row, 10).Value = -0.29          # P&L $
    $ws.Cells.Item($row, 11).Value = 99.81999999999999  # Capital After
    $ws.Cells.Item($row, 12).Value = 0              # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0              # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6            # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"   # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.14           # Duration (min)
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
